$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 334 and 335, pushing the existing rows
# (old 334..351) down to become 336..353.
$ws.Rows("334:335").Insert()

# --- New row 334: Lane Late / Provincia de Quillota -----------------------
$ws.Range("A334").Value = 5
$ws.Range("B334").Value = "Macroferia Regional de Talca"
$ws.Range("C334").Value = "Maule"
$ws.Range("D334").Value = 44516
$ws.Range("E334").Value = 7
$ws.Range("F334").Value = "Fruta"
$ws.Range("G334").Value = 100102
$ws.Range("H334").Value = "Cítricos"
$ws.Range("I334").Value = 100102005
$ws.Range("J334").Value = "Naranja"
$ws.Range("K334").Value = "Lane Late"
$ws.Range("L334").Value = "Primera"
$ws.Range("M334").Value = 350
$ws.Range("N334").Value = 8000
$ws.Range("O334").Value = 8000
$ws.Range("P334").Value = 8000
$ws.Range("Q334").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R334").Value = "Provincia de Quillota"
$ws.Range("S334").Value = 533
$ws.Range("T334").Value = 15

# --- New row 335: Navel Late / Provincia de Quillota -----------------------
$ws.Range("A335").Value = 5
$ws.Range("B335").Value = "Macroferia Regional de Talca"
$ws.Range("C335").Value = "Maule"
$ws.Range("D335").Value = 44516
$ws.Range("E335").Value = 7
$ws.Range("F335").Value = "Fruta"
$ws.Range("G335").Value = 100102
$ws.Range("H335").Value = "Cítricos"
$ws.Range("I335").Value = 100102005
$ws.Range("J335").Value = "Naranja"
$ws.Range("K335").Value = "Navel Late"
$ws.Range("L335").Value = "Primera"
$ws.Range("M335").Value = 360
$ws.Range("N335").Value = 8000
$ws.Range("O335").Value = 8000
$ws.Range("P335").Value = 8000
$ws.Range("Q335").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R335").Value = "Provincia de Quillota"
$ws.Range("S335").Value = 533
$ws.Range("T335").Value = 15
